$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 42: add missing "属性加成" ("/") value in column C ---
$ws.Range("C42").Value = "/"

# --- New data rows 44-49 ---
# Values are entered in an order that mirrors how the shared-string table
# grows (column-by-column across the related rows), matching the source
# workbook's authoring pattern.

# Row 44 & 45 (料理大赛 related rows)
$ws.Range("A44").Value = "—团团圆圆—"
$ws.Range("A45").Value = "—料理达人—"

$ws.Range("D44").Value = "1天"
$ws.Range("D45").Value = "1天"

$ws.Range("E44").Value = "2024.02.22-2024.03.06"
$ws.Range("E45").Value = "2024.02.22-2024.03.06"

$ws.Range("F44").Value = "料理大赛积分榜第一名"
$ws.Range("F45").Value = "料理大赛积分榜第二名至第十名"

$ws.Range("G44").Value = "料理大赛"
$ws.Range("G45").Value = "料理大赛"

$ws.Range("B44").Value = "金色"
$ws.Range("C44").Value = "/"
$ws.Range("B45").Value = "紫色"
$ws.Range("C45").Value = "/"

# Row 46 (节日连线 - 学业有成连连看)
$ws.Range("A46").Value = "—头好痒要长脑子了—"
$ws.Range("B46").Value = "紫色+底纹"
$ws.Range("G46").Value = "节日连线-学业有成连连看"
$ws.Range("E46").Value = "2024.02.14-2024.03.07"
$ws.Range("C46").Value = "/"
$ws.Range("D46").Value = "永久"

# Row 47 & 48 (全国护肝日 related rows)
$ws.Range("A47").Value = "—躺一下怎么了—"
$ws.Range("A48").Value = "—躺一下怎么了—"

$ws.Range("E47").Value = "2024.03.14-2024.03.21"
$ws.Range("E48").Value = "2024.03.14-2024.03.21"

$ws.Range("F47").Value = "完成护肝指南全部任务"
$ws.Range("F48").Value = "完成护肝指南任一任务"

$ws.Range("G47").Value = "全国护肝日"
$ws.Range("G48").Value = "全国护肝日"

$ws.Range("B47").Value = "紫色+底纹"
$ws.Range("C47").Value = "/"
$ws.Range("D47").Value = "永久"
$ws.Range("B48").Value = "蓝色"
$ws.Range("C48").Value = "/"
$ws.Range("D48").Value = "永久"

# Row 49 (樱花绽放)
$ws.Range("A49").Value = "—春日樱花—"
$ws.Range("E49").Value = "2024.03.21-2024.04.03"
$ws.Range("G49").Value = "樱花绽放"
$ws.Range("F49").Value = "樱花绽放打卡4个地方"
$ws.Range("B49").Value = "紫色+底纹"
$ws.Range("D49").Value = "永久"

# --- Apply matching font styling to column A of new rows (copy formats only) ---
$ws.Range("A3").Copy()
$ws.Range("A44").PasteSpecial(-4122)

$ws.Range("A4").Copy()
$ws.Range("A45").PasteSpecial(-4122)
$ws.Range("A46").PasteSpecial(-4122)
$ws.Range("A49").PasteSpecial(-4122)

$ws.Range("A38").Copy()
$ws.Range("A47").PasteSpecial(-4122)

$ws.Range("A34").Copy()
$ws.Range("A48").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Column G width update ---
$ws.Columns.Item(7).ColumnWidth = 21.5

# --- Update sheet view: select B50 (final cursor position) ---
$ws.Range("A45").Select()
$ws.Range("B50").Select()
